# Remove the third data row (A3:D3) from Sheet1, matching the "delete
# selected row contents, keep formatting" action: the row's cell contents
# (and the shared-string entries that become unused) disappear, but the
# styled cell C3 keeps its style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("A3:D3")

# Clear the values/formulas but keep cell formatting (C3 keeps its style).
$target.ClearContents()

# Leave the range selected, like it was right after the Delete keypress.
$target.Select()
